$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RiskRegister")
$ws.Range("G2").Value = "4. Significant Impact"
